# Append three new call-log rows (rows 4-6) to the Incomplete_Calls sheet,
# matching the "enhnacements in call_queue monitoring" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Vanshika panjwani - minimal_interaction, duration 0
$ws.Cells.Item(4, 1).Value = "Vanshika panjwani"
$ws.Cells.Item(4, 2).Value = "'917823844614"
$ws.Cells.Item(4, 3).Value = "24 MG Road, Bengaluru"
$ws.Cells.Item(4, 4).Value = "'28"
$ws.Cells.Item(4, 5).Value = "Male"
$ws.Cells.Item(4, 6).Value = "2025-06-25 19:35:46"
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = "minimal_interaction"
$ws.Cells.Item(4, 9).Value = "Very few exchanges in conversation"

# Row 5: Vanshika panjwani - minimal_interaction, duration 88
$ws.Cells.Item(5, 1).Value = "Vanshika panjwani"
$ws.Cells.Item(5, 2).Value = "'917823844614"
$ws.Cells.Item(5, 3).Value = "24 MG Road, Bengaluru"
$ws.Cells.Item(5, 4).Value = "'28"
$ws.Cells.Item(5, 5).Value = "Male"
$ws.Cells.Item(5, 6).Value = "2025-06-25 20:32:23"
$ws.Cells.Item(5, 7).Value = 88
$ws.Cells.Item(5, 8).Value = "minimal_interaction"
$ws.Cells.Item(5, 9).Value = "Very few exchanges in conversation"

# Row 6: Vanshika panjwani - minimal_interaction, duration 1
$ws.Cells.Item(6, 1).Value = "Vanshika panjwani"
$ws.Cells.Item(6, 2).Value = "'917823844614"
$ws.Cells.Item(6, 3).Value = "24 MG Road, Bengaluru"
$ws.Cells.Item(6, 4).Value = "'28"
$ws.Cells.Item(6, 5).Value = "Male"
$ws.Cells.Item(6, 6).Value = "2025-06-25 20:34:46"
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = "minimal_interaction"
$ws.Cells.Item(6, 9).Value = "Very few exchanges in conversation"
